$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ST.07")
$ws.Range("E10:J24").Value = 0
